# Atualização automática via Streamlit (13/11/2025 19:09)
# Appends a new data row (row 4) to the PRINCIPAL sheet of SALDO_PECAS.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "AM"
$ws.Range("B4").Value = "00x0098"

# SUB1/SUB2/SUB3 stay blank for this entry, but the cells themselves must
# still exist in the sheet (same shape as rows 2 and 3), so touch each one
# (a no-op format nudge at the default size) to materialise it without
# leaving the value empty-string assignment from wiping the cell back out.
$ws.Range("C4").Font.Size = 11
$ws.Range("D4").Font.Size = 11
$ws.Range("E4").Font.Size = 11

$ws.Range("F4").Value = "p"
$ws.Range("G4").Value = "p"
$ws.Range("H4").Value = "p - (p 01/08/25_24h) - AM"

# DATA_FIM / DATA_VERIFICACAO look like dates ("dd/mm/yy") but must be kept
# as literal text, matching the source data (inline string "01/08/25"),
# not converted to a date serial number.
$ws.Range("I4").NumberFormat = "@"
$ws.Range("I4").Value = "01/08/25"

$ws.Range("J4").Value = "24h"

$ws.Range("K4").NumberFormat = "@"
$ws.Range("K4").Value = "13/11/25"

$ws.Range("L4").Value = "DENTRO"
